$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 97.5
$ws.Range("I6").Value = 96.666664
$ws.Range("J6").Value = 100
$ws.Range("K6").Value = 289.999992
$ws.Range("L6").Value = 300
$ws.Range("M6").Value = -177.999992
$ws.Range("N6").Value = -524
$ws.Range("H12").Value = 181.57143
$ws.Range("I12").Value = 149.6
$ws.Range("K12").Value = 149.6
$ws.Range("M12").Value = 20.40000000000001
$ws.Range("H20").Value = 1499.5
$ws.Range("I20").Value = 999
$ws.Range("K20").Value = 999
$ws.Range("M20").Value = -769
$ws.Range("H28").Value = 1075.4
$ws.Range("I28").Value = 844.25
$ws.Range("J28").Value = 2000
$ws.Range("K28").Value = 844.25
$ws.Range("L28").Value = 2000
$ws.Range("M28").Value = -359.25
$ws.Range("N28").Value = -2970
$ws.Range("H35").Value = 1499.5
$ws.Range("I35").Value = 999
$ws.Range("K35").Value = 999
$ws.Range("M35").Value = -620
$ws.Range("H55").Value = 656.6667
$ws.Range("I55").Value = 740
$ws.Range("K55").Value = 740
$ws.Range("M55").Value = -526
$ws.Range("H64").Value = 2166.6667
$ws.Range("I64").Value = 2000
$ws.Range("J64").Value = 2500
$ws.Range("K64").Value = 2000
$ws.Range("L64").Value = 2500
$ws.Range("M64").Value = -1752
$ws.Range("N64").Value = -2996
$ws.Range("H67").Value = 2166.6667
$ws.Range("I67").Value = 2000
$ws.Range("J67").Value = 2500
$ws.Range("K67").Value = 2000
$ws.Range("L67").Value = 2500
$ws.Range("M67").Value = -1142
$ws.Range("N67").Value = -4216
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 90.333336
$ws.Range("I5").Value = 122.5
$ws.Range("J5").Value = 26
$ws.Range("K5").Value = 122.5
$ws.Range("L5").Value = 26
$ws.Range("M5").Value = -10.5
$ws.Range("N5").Value = -250
$ws.Range("H37").Value = 10000
$ws.Range("I37").Value = 10000
$ws.Range("K37").Value = 10000
$ws.Range("M37").Value = -9727
$ws.Range("H61").Value = 2938.5454
$ws.Range("I61").Value = 2190.1428
$ws.Range("K61").Value = 2190.1428
$ws.Range("M61").Value = -1978.1428
$ws.Range("H95").Value = 44221.6
$ws.Range("J95").Value = 44221.6
$ws.Range("L95").Value = 44221.6
$ws.Range("N95").Value = -49713.6
$ws.Range("H132").Value = 5333.857
$ws.Range("I132").Value = 3467.4
$ws.Range("K132").Value = 10402.2
$ws.Range("M132").Value = -7872.200000000001
$ws.Range("H136").Value = 2938.5454
$ws.Range("I136").Value = 2190.1428
$ws.Range("K136").Value = 6570.428400000001
$ws.Range("M136").Value = -4020.428400000001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 90.333336
$ws.Range("I4").Value = 122.5
$ws.Range("J4").Value = 26
$ws.Range("K4").Value = 122.5
$ws.Range("L4").Value = 26
$ws.Range("M4").Value = -7.5
$ws.Range("N4").Value = -256
$ws.Range("H64").Value = 1182.5
$ws.Range("I64").Value = 990
$ws.Range("J64").Value = 1375
$ws.Range("K64").Value = 990
$ws.Range("L64").Value = 1375
$ws.Range("M64").Value = -765
$ws.Range("N64").Value = -1825
$ws.Range("H67").Value = 1182.5
$ws.Range("I67").Value = 990
$ws.Range("J67").Value = 1375
$ws.Range("K67").Value = 990
$ws.Range("L67").Value = 1375
$ws.Range("M67").Value = -210
$ws.Range("N67").Value = -2935
$ws.Range("H82").Value = 16999
$ws.Range("I82").Value = 16999
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 16999
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -16616
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 16999
$ws.Range("I85").Value = 16999
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 16999
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -15673
$ws.Range("N85").ClearContents()
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("M97").ClearContents()
$ws.Range("H134").Value = 1498.75
$ws.Range("I134").Value = 1498.75
$ws.Range("K134").Value = 4496.25
$ws.Range("M134").Value = -1961.25

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4373.1665
$ws.Range("I58").Value = 2503
$ws.Range("J58").Value = 5308.25
$ws.Range("K58").Value = 2503
$ws.Range("L58").Value = 5308.25
$ws.Range("M58").Value = -2300
$ws.Range("N58").Value = -5714.25
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H136").Value = 4373.1665
$ws.Range("I136").Value = 2503
$ws.Range("J136").Value = 5308.25
$ws.Range("K136").Value = 7509
$ws.Range("L136").Value = 15924.75
$ws.Range("M136").Value = -4959
$ws.Range("N136").Value = -21024.75

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 3666.6667
$ws.Range("I55").Value = 2500
$ws.Range("J55").Value = 4250
$ws.Range("K55").Value = 7500
$ws.Range("L55").Value = 12750
$ws.Range("M55").Value = -7323
$ws.Range("N55").Value = -13104
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H68").Value = 802
$ws.Range("J68").Value = 802
$ws.Range("L68").Value = 2406
$ws.Range("N68").Value = -4028
$ws.Range("H71").Value = 802
$ws.Range("J71").Value = 802
$ws.Range("L71").Value = 7218
$ws.Range("N71").Value = -15330
$ws.Range("H131").Value = 2974.75
$ws.Range("J131").Value = 900
$ws.Range("L131").Value = 2700
$ws.Range("N131").Value = -12780

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 33334
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 33334
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 33334
$ws.Range("M10").ClearContents()
$ws.Range("N10").Value = -33672
$ws.Range("H15").Value = 30000
$ws.Range("J15").Value = 30000
$ws.Range("L15").Value = 30000
$ws.Range("N15").Value = -30576
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()
$ws.Range("H24").Value = 13169.333
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 13169.333
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 13169.333
$ws.Range("M24").ClearContents()
$ws.Range("N24").Value = -13515.333
$ws.Range("H30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("N30").ClearContents()
$ws.Range("H49").Value = 39300
$ws.Range("J49").Value = 39300
$ws.Range("L49").Value = 39300
$ws.Range("N49").Value = -39668
$ws.Range("H81").Value = 30000
$ws.Range("J81").Value = 30000
$ws.Range("L81").Value = 30000
$ws.Range("N81").Value = -31996
$ws.Range("H84").Value = 30000
$ws.Range("J84").Value = 30000
$ws.Range("L84").Value = 90000
$ws.Range("N84").Value = -99984

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 4300
$ws.Range("I16").Value = 4266.6665
$ws.Range("K16").Value = 4266.6665
$ws.Range("M16").Value = -4096.6665
$ws.Range("H22").Value = 750.5
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 750.5
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 750.5
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -1340.5
$ws.Range("H27").Value = 750.5
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 750.5
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 750.5
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -964.5
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("H68").Value = 1350
$ws.Range("I68").Value = 700
$ws.Range("J68").Value = 2000
$ws.Range("K68").Value = 700
$ws.Range("L68").Value = 2000
$ws.Range("M68").Value = 49
$ws.Range("N68").Value = -3498
$ws.Range("H71").Value = 1350
$ws.Range("I71").Value = 700
$ws.Range("J71").Value = 2000
$ws.Range("K71").Value = 3500
$ws.Range("L71").Value = 10000
$ws.Range("M71").Value = 244
$ws.Range("N71").Value = -17488

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 900000
$ws.Range("I31").Value = 900000
$ws.Range("K31").Value = 900000
$ws.Range("M31").Value = -899652
$ws.Range("H64").Value = 50000
$ws.Range("I64").Value = 50000
$ws.Range("K64").Value = 50000
$ws.Range("M64").Value = -49752
$ws.Range("H67").Value = 50000
$ws.Range("I67").Value = 50000
$ws.Range("K67").Value = 50000
$ws.Range("M67").Value = -49142
$ws.Range("H70").Value = 39368.332
$ws.Range("I70").Value = 38000
$ws.Range("K70").Value = 38000
$ws.Range("M70").Value = -37685
$ws.Range("H73").Value = 39368.332
$ws.Range("I73").Value = 38000
$ws.Range("K73").Value = 38000
$ws.Range("M73").Value = -37685
